$d = $word.ActiveDocument

# Find the paragraph that holds the knitr warning source-code block and
# remove it entirely (including its paragraph mark), leaving the
# "2019-04-17" date paragraph directly followed by the "Introduction"
# heading paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Warning: package*knitr*built under R version 3.4.4*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.Delete()
}
